$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update answer text for the two existing answer rows
$ws.Range("B2").Value = "Bob"
$ws.Range("B3").Value = "Kari"

# Widen column B to fit the new content (matches the bestFit width recorded after edit)
$ws.Range("B2:B3").Columns.AutoFit()
$ws.Columns.Item(2).ColumnWidth = 26.140625

# Move the active selection to B8
$ws.Range("B8").Select()
